$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand header row (N1:P1) by copying the style from M1, then set the date labels
$ws.Range("M1").Copy($ws.Range("N1:P1"))
$ws.Cells.Item(1, 14).Value = "31/12/2023"
$ws.Cells.Item(1, 15).Value = "31/03/2024"
$ws.Cells.Item(1, 16).Value = "30/06/2024"

# Fill in the new quarterly figures (31/12/2023, 31/03/2024, 30/06/2024) for each data row
$ws.Cells.Item(2, 14).Value = 10186692.608
$ws.Cells.Item(2, 15).Value = 11096087.552
$ws.Cells.Item(2, 16).Value = 10525810.688
$ws.Cells.Item(3, 14).Value = 5974663.168
$ws.Cells.Item(3, 15).Value = 6853248
$ws.Cells.Item(3, 16).Value = 6223273.984
$ws.Cells.Item(4, 14).Value = 472702.016
$ws.Cells.Item(4, 15).Value = 1137089.024
$ws.Cells.Item(4, 16).Value = 720732.032
$ws.Cells.Item(5, 14).Value = 528792
$ws.Cells.Item(5, 15).Value = 605033.9840000001
$ws.Cells.Item(5, 16).Value = 735886.976
$ws.Cells.Item(6, 14).Value = 2453331.968
$ws.Cells.Item(6, 15).Value = 2663852.032
$ws.Cells.Item(6, 16).Value = 2348581.888
$ws.Cells.Item(7, 14).Value = 1980075.008
$ws.Cells.Item(7, 15).Value = 1901744
$ws.Cells.Item(7, 16).Value = 1910578.944
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(9, 14).Value = 386492.992
$ws.Cells.Item(9, 15).Value = 381536.992
$ws.Cells.Item(9, 16).Value = 363336
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(11, 14).Value = 153268.992
$ws.Cells.Item(11, 15).Value = 163992
$ws.Cells.Item(11, 16).Value = 144158
$ws.Cells.Item(12, 14).Value = 1010518.976
$ws.Cells.Item(12, 15).Value = 1071468.992
$ws.Cells.Item(12, 16).Value = 1148272
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(16, 14).Value = 7925
$ws.Cells.Item(16, 15).Value = 7410
$ws.Cells.Item(16, 16).Value = 20466
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(19, 14).Value = 544638.976
$ws.Cells.Item(19, 15).Value = 600195.008
$ws.Cells.Item(19, 16).Value = 641819.008
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(21, 15).Value = 0
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(22, 14).Value = 13045
$ws.Cells.Item(22, 15).Value = 12480
$ws.Cells.Item(22, 16).Value = 11977
$ws.Cells.Item(23, 14).Value = 554435.008
$ws.Cells.Item(23, 15).Value = 549803.008
$ws.Cells.Item(23, 16).Value = 547462.0159999999
$ws.Cells.Item(24, 14).Value = 2634031.104
$ws.Cells.Item(24, 15).Value = 2609088
$ws.Cells.Item(24, 16).Value = 2594825.984
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(25, 15).Value = 0
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(26, 14).Value = 10186692.608
$ws.Cells.Item(26, 15).Value = 11096087.552
$ws.Cells.Item(26, 16).Value = 10525810.688
$ws.Cells.Item(27, 14).Value = 3056590.08
$ws.Cells.Item(27, 15).Value = 2953587.968
$ws.Cells.Item(27, 16).Value = 2622516.992
$ws.Cells.Item(28, 14).Value = 108718
$ws.Cells.Item(28, 15).Value = 106004
$ws.Cells.Item(28, 16).Value = 112677
$ws.Cells.Item(29, 14).Value = 1716115.968
$ws.Cells.Item(29, 15).Value = 1642434.944
$ws.Cells.Item(29, 16).Value = 1583928.96
$ws.Cells.Item(30, 14).Value = 74344
$ws.Cells.Item(30, 15).Value = 89285
$ws.Cells.Item(30, 16).Value = 67840
$ws.Cells.Item(31, 14).Value = 391745.984
$ws.Cells.Item(31, 15).Value = 410950.016
$ws.Cells.Item(31, 16).Value = 237732
$ws.Cells.Item(32, 14).Value = 0
$ws.Cells.Item(32, 15).Value = 0
$ws.Cells.Item(32, 16).Value = 0
$ws.Cells.Item(33, 14).Value = 106138
$ws.Cells.Item(33, 15).Value = 102576
$ws.Cells.Item(33, 16).Value = 43482
$ws.Cells.Item(34, 14).Value = 659527.9840000001
$ws.Cells.Item(34, 15).Value = 602337.9840000001
$ws.Cells.Item(34, 16).Value = 576857.008
$ws.Cells.Item(35, 14).Value = 0
$ws.Cells.Item(35, 15).Value = 0
$ws.Cells.Item(35, 16).Value = 0
$ws.Cells.Item(36, 14).Value = 0
$ws.Cells.Item(36, 15).Value = 0
$ws.Cells.Item(36, 16).Value = 0
$ws.Cells.Item(37, 14).Value = 3840742.912
$ws.Cells.Item(37, 15).Value = 4856220.16
$ws.Cells.Item(37, 16).Value = 4644007.936
$ws.Cells.Item(38, 14).Value = 2835177.984
$ws.Cells.Item(38, 15).Value = 3845082.88
$ws.Cells.Item(38, 16).Value = 3665315.072
$ws.Cells.Item(39, 14).Value = 0
$ws.Cells.Item(39, 15).Value = 0
$ws.Cells.Item(39, 16).Value = 0
$ws.Cells.Item(40, 14).Value = 888241.9840000001
$ws.Cells.Item(40, 15).Value = 893564.992
$ws.Cells.Item(40, 16).Value = 850033.9840000001
$ws.Cells.Item(41, 14).Value = 50266
$ws.Cells.Item(41, 15).Value = 39943
$ws.Cells.Item(41, 16).Value = 32451
$ws.Cells.Item(42, 14).Value = 0
$ws.Cells.Item(42, 15).Value = 0
$ws.Cells.Item(42, 16).Value = 0
$ws.Cells.Item(43, 14).Value = 67057
$ws.Cells.Item(43, 15).Value = 77629
$ws.Cells.Item(43, 16).Value = 96208
$ws.Cells.Item(44, 14).Value = 0
$ws.Cells.Item(44, 15).Value = 0
$ws.Cells.Item(44, 16).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 0
$ws.Cells.Item(45, 16).Value = 0
$ws.Cells.Item(46, 14).Value = 0
$ws.Cells.Item(46, 15).Value = 0
$ws.Cells.Item(46, 16).Value = 0
$ws.Cells.Item(47, 14).Value = 3289359.872
$ws.Cells.Item(47, 15).Value = 3286279.936
$ws.Cells.Item(47, 16).Value = 3259286.016
$ws.Cells.Item(48, 14).Value = 2549391.872
$ws.Cells.Item(48, 15).Value = 2549391.872
$ws.Cells.Item(48, 16).Value = 2549391.872
$ws.Cells.Item(49, 14).Value = -268287.008
$ws.Cells.Item(49, 15).Value = -277972
$ws.Cells.Item(49, 16).Value = -276344
$ws.Cells.Item(50, 14).Value = 0
$ws.Cells.Item(50, 15).Value = 0
$ws.Cells.Item(50, 16).Value = 0
$ws.Cells.Item(51, 14).Value = 1008254.976
$ws.Cells.Item(51, 15).Value = 1014860.032
$ws.Cells.Item(51, 16).Value = 986238.0159999999
$ws.Cells.Item(52, 14).Value = 0
$ws.Cells.Item(52, 15).Value = 0
$ws.Cells.Item(52, 16).Value = 0
$ws.Cells.Item(53, 14).Value = 0
$ws.Cells.Item(53, 15).Value = 0
$ws.Cells.Item(53, 16).Value = 0
$ws.Cells.Item(54, 14).Value = 0
$ws.Cells.Item(54, 15).Value = 0
$ws.Cells.Item(54, 16).Value = 0
$ws.Cells.Item(55, 14).Value = 0
$ws.Cells.Item(55, 15).Value = 0
$ws.Cells.Item(55, 16).Value = 0
$ws.Cells.Item(56, 14).Value = 0
$ws.Cells.Item(56, 15).Value = 0
$ws.Cells.Item(56, 16).Value = 0
$ws.Cells.Item(59, 14).Value = 2905142.272
$ws.Cells.Item(59, 15).Value = 2958138.88
$ws.Cells.Item(59, 16).Value = 2753534.976
$ws.Cells.Item(60, 14).Value = -2484379.904
$ws.Cells.Item(60, 15).Value = -2559045.888
$ws.Cells.Item(60, 16).Value = -2356251.904
$ws.Cells.Item(61, 14).Value = 420762.048
$ws.Cells.Item(61, 15).Value = 399092.992
$ws.Cells.Item(61, 16).Value = 397283.008
$ws.Cells.Item(62, 14).Value = -85382.008
$ws.Cells.Item(62, 15).Value = -98012
$ws.Cells.Item(62, 16).Value = -102179
$ws.Cells.Item(63, 14).Value = -226992.992
$ws.Cells.Item(63, 15).Value = -236132
$ws.Cells.Item(63, 16).Value = -238956.992
$ws.Cells.Item(64, 14).Value = -1477
$ws.Cells.Item(64, 15).Value = 3199
$ws.Cells.Item(64, 16).Value = -2517
$ws.Cells.Item(65, 14).Value = 16330
$ws.Cells.Item(65, 15).Value = 14789
$ws.Cells.Item(65, 16).Value = 3897
$ws.Cells.Item(66, 14).Value = -4296
$ws.Cells.Item(66, 15).Value = -12713
$ws.Cells.Item(66, 16).Value = -23078
$ws.Cells.Item(67, 14).Value = -475
$ws.Cells.Item(67, 15).Value = -467
$ws.Cells.Item(67, 16).Value = -403
$ws.Cells.Item(68, 14).Value = -88132
$ws.Cells.Item(68, 15).Value = -112770
$ws.Cells.Item(68, 16).Value = -165968
$ws.Cells.Item(69, 14).Value = 62172.992
$ws.Cells.Item(69, 15).Value = 30536
$ws.Cells.Item(69, 16).Value = 39430
$ws.Cells.Item(70, 14).Value = -150304.992
$ws.Cells.Item(70, 15).Value = -143306
$ws.Cells.Item(70, 16).Value = -205398
$ws.Cells.Item(74, 14).Value = 30336.992
$ws.Cells.Item(74, 15).Value = -43013
$ws.Cells.Item(74, 16).Value = -131922
$ws.Cells.Item(75, 14).Value = 74370
$ws.Cells.Item(75, 15).Value = -15691
$ws.Cells.Item(75, 16).Value = -5376
$ws.Cells.Item(76, 14).Value = 156536
$ws.Cells.Item(76, 15).Value = 61747
$ws.Cells.Item(76, 16).Value = 49581
$ws.Cells.Item(79, 14).Value = 0
$ws.Cells.Item(79, 15).Value = 0
$ws.Cells.Item(79, 16).Value = 0
$ws.Cells.Item(80, 14).Value = 261243.008
$ws.Cells.Item(80, 15).Value = 3043
$ws.Cells.Item(80, 16).Value = -87717

# Blank separator rows: replicate the empty inline-string cells used by columns B:M
$ws.Range("M57:M58").Copy($ws.Range("N57:P58"))
$ws.Range("M71:M73").Copy($ws.Range("N71:P73"))
$ws.Range("M77:M78").Copy($ws.Range("N77:P78"))
